# Apply the commit "Updated notebook, reran simulation":
#  - Two new measurement cases ("Holden" and "Rizzie Spiral") were added to the
#    dataset, inserted right after "Spiral5" / before "RotRing OmegaMax-90".
#    This shifts every subsequent row down by two positions.
#  - One existing case was renamed: "Thomas Hex" -> "Matthies Hex".
#  - The simulation was rerun, which is why numbers "shifted" along with rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 4 (pushes old rows 4..29 down to rows 6..31,
# and the former last two rows 28/29 naturally slide down into new rows 30/31).
$ws.Range("A4:A5").EntireRow.Insert()

# Copy the "index column" formatting (bold, border, centered) from row 6 down
# onto the two freshly-inserted rows, matching the rest of column A.
$ws.Range("A6").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)

# --- Fill in the new "Holden" row (row 4) ---
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "Holden"

$row4 = @{
  3  = 0.9950484231798352
  4  = 0.9950484231798352
  5  = 0.9177275512883423
  6  = 1.091612077757663
  7  = 0.9460216726130078
  8  = 0.7892930973116186
  9  = 0.6979339306661129
  10 = 0.9503055010233554
  11 = 4.723790133246625
  12 = 0.6979339306661129
  13 = 0.9950484231797458
  14 = 0.9950484231798352
  15 = 4.723790133246625
  16 = 2.710862031956369
  17 = 2.820758842267483
  18 = 2.138924162364191
  19 = 2.11315053840036
  20 = 2.138924162364191
  21 = 1.833625009595229
  22 = 1.66590969231215
  23 = 1.38896654838582
}
foreach ($col in $row4.Keys) {
  $ws.Cells.Item(4, $col).Value = $row4[$col]
}

# --- Fill in the new "Rizzie Spiral" row (row 5) ---
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "Rizzie Spiral"

$row5 = @{
  3  = 11.85567036476542
  4  = 11.85567036476542
  5  = 3.960888691994206
  6  = 3.604885020557555
  7  = 3.957814197546023
  8  = 3.878833375896534
  9  = 8.979342344833363
  10 = -0.001134647219809994
  11 = 0.6933168727090753
  12 = 8.979342344833363
  13 = 11.85567036476542
  14 = 11.85567036476542
  15 = 0.6933168727090753
  16 = 4.836329608771219
  17 = 2.327102782351641
  18 = 7.176109860769285
  19 = 4.544515969845548
  20 = 7.176109860769285
  21 = 6.372304568575515
  22 = 7.468977727813495
  23 = 4.616202027635295
}
foreach ($col in $row5.Keys) {
  $ws.Cells.Item(5, $col).Value = $row5[$col]
}

# --- Rename "Thomas Hex" -> "Matthies Hex" (now located at row 11 after the shift) ---
$ws.Cells.Item(11,2).Value = "Matthies Hex"
